$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Remove the two rows whose entire content was dropped ---
# Old Sr.No 9 "Policy ID Validation" (row 11)
$ws.Rows(11).Delete()
# Old Sr.No 14 "Partner Certificate Signing and RE-issueing" (now row 15 after the first delete)
$ws.Rows(15).Delete()

# --- Step 2: Update cell text/content that changed wording or was replaced ---
$ws.Range("E5").Value = "1. Validate length of a License Key as configured and respond as mentioned below`na. If found valid, respond with ""VALID""`nb. if found invalid, respond with ""INVALID"""
$ws.Range("E9").Value = "1. Validate length of a Partner ID as configured and respond as mentioned below`na. If found valid, respond with ""VALID""`nb. if found invalid, respond with ""INVALID"""
$ws.Range("D10").Value = "Map Policies to Partners"
$ws.Range("E10").Value = "1. Map following Policies to Partners`na. Auth Policies ( can be Mandatory/Non-Mandatory)`n     1. OTP Trigger `n     2. OTP Authentication`n     3. Demo Authentication `n     4. Biometric Authentication - FMR Data Match `n     5. Biometric Authentication - IIR Data Match  `n     6. Biometric Authentication - FID Data Match `nb. E-Kyc Policies (can be Required/Not Required)`n    1. eKYC - all combinations of eKYC demo fields "
$ws.Range("D11").Value = "Retrieve Policies based on Partner ID"
$ws.Range("E11").Value = "1. Receive request to retreive policies based on Partner ID`n2. Respond appropirately if Partner ID does not exist"
$ws.Range("E12").Value = "1. Receive request to register a Partner with follwing parameters`na. Partner Name`nb. Partner Contact Name`nc. Partner Phone`nd. Partner Email ID`n2. Issue and Map Partner ID`n3. Map Policies to the Partner`na. Multiple Policies can be mapped to a Partner`nb. A Partner can have a policy for both Auth and E-KYC`n4. Store the Partner in MOSIP"
$ws.Range("E13").Value = "1. Receive a request to map MISP to a Partner with MISP ID and Partner ID as Input`n2. There can be a many-to-mapping between MISPs and Partners"
$ws.Range("D14").Value = "Partner Certificate Validation"
$ws.Range("E14").Value = "1. Upload Digital Certificate on Admin Portal for a Partner`n2. Verify CA Authority of the certificate`n3. Sign the certificate with MOSIP Certificate`n4. Respond to the source with the re-issued certificate`n5. Certificate will be uploaded by the MOSIP admin. The Registered Partner will send the certificate to the MOSIP Admin through ofline process. Re-issued certificate will be sent to the Partner by MOSIP admin through notification/offline process`n6. Private key to change priodically as per the Key Rotation Policy set by admin"
$ws.Range("E15").Value = "1. Distribute Public Key to Partners for encrypting the Auth Request befoe sending it to the MOSIP`n2. Public key needs to be distributed priodically whenever the correspinding Private Key is rotated"

# --- Step 3: Row heights that changed because of the re-wrapped text ---
$ws.Rows(10).RowHeight = 145
$ws.Rows(11).RowHeight = 29
$ws.Rows(14).RowHeight = 145

# --- Step 4: Restore cursor / selection position shown in the saved file ---
$ws.Range("E15").Select()
